$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for the new "Código" column
$ws.Columns.Item(1).Insert()

# Set the header for the new first column
$ws.Range("A1").Value = "Código"

# Update the active selection to match the authored change
$ws.Range("B4").Select()
